$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Bump the "datetimeFigureOut" auto-date placeholders (slide master, every
#    slide layout, and the notes master) from 22/06/2020 to 08/07/2020.
# ---------------------------------------------------------------------------
$oldDate = "22/06/2020"
$newDate = "08/07/2020"

function Update-DateShape($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        if ($shp.HasTextFrame -and $shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$sm = $p.SlideMaster
Update-DateShape $sm.Shapes

for ($i = 1; $i -le $sm.CustomLayouts.Count; $i++) {
    $cl = $sm.CustomLayouts.Item($i)
    Update-DateShape $cl.Shapes
}

$nm = $p.NotesMaster
Update-DateShape $nm.Shapes

# ---------------------------------------------------------------------------
# 2) Slide 2: bump the credits textbox from "v1.0.1" to "v1.0.2".
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $shp = $s2.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        $full = $tr.Text
        $idx = $full.IndexOf("Outcome game (v1.0.1) ")
        if ($idx -ge 0) {
            $sub = $tr.Characters($idx + 1, 22)
            $sub.Text = "Outcome game (v1.0.2) "
        }
    }
}

# ---------------------------------------------------------------------------
# 3) Slide 57: recenter the "Don't measure Output. Measure Outcomes" textbox
#    (move it up and make it taller, keeping left/width unchanged).
# ---------------------------------------------------------------------------
$s57 = $p.Slides.Item(57)
for ($i = 1; $i -le $s57.Shapes.Count; $i++) {
    $shp = $s57.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 5") {
        $shp.Top = 115.23590551181101
        $shp.Height = 189.02812523622047
    }
}
